# Recomputed NATMI TPM-derived ligand-receptor metrics for Sheet1 (App-Fpr2).
# Every cell below is a plain numeric overwrite; no rows/cols/styles/strings change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 115.2213693333333
$ws.Range("H2").Value = 345.664108
$ws.Range("I2").Value = 0.2787408744545015
$ws.Range("J2").Value = 0.2787408744545015
$ws.Range("M2").Value = 0.05237733333333333
$ws.Range("N2").Value = 0.157132
$ws.Range("O2").Value = 0.01249705432598326
$ws.Range("P2").Value = 0.01249705432598326
$ws.Range("Q2").Value = 6.03498806869511
$ws.Range("R2").Value = 54.314892618256
$ws.Range("S2").Value = 0.003483439850929985
$ws.Range("T2").Value = 0.003483439850929983

# Row 3
$ws.Range("G3").Value = 115.2213693333333
$ws.Range("H3").Value = 345.664108
$ws.Range("I3").Value = 0.2787408744545015
$ws.Range("J3").Value = 0.2787408744545015
$ws.Range("O3").Value = 0.1346970137168397
$ws.Range("P3").Value = 0.1346970137168397
$ws.Range("Q3").Value = 65.04691821494755
$ws.Range("R3").Value = 585.422263934528
$ws.Range("S3").Value = 0.03754556338984189
$ws.Range("T3").Value = 0.03754556338984189

# Row 4
$ws.Range("G4").Value = 115.2213693333333
$ws.Range("H4").Value = 345.664108
$ws.Range("I4").Value = 0.2787408744545015
$ws.Range("J4").Value = 0.2787408744545015
$ws.Range("O4").Value = 0.020817156814363
$ws.Range("P4").Value = 0.020817156814363
$ws.Range("Q4").Value = 10.05287243871778
$ws.Range("R4").Value = 90.47585194846
$ws.Range("S4").Value = 0.005802592494092029
$ws.Range("T4").Value = 0.005802592494092027

# Row 5
$ws.Range("G5").Value = 115.2213693333333
$ws.Range("H5").Value = 345.664108
$ws.Range("I5").Value = 0.2787408744545015
$ws.Range("J5").Value = 0.2787408744545015
$ws.Range("M5").Value = 3.48701
$ws.Range("N5").Value = 10.46103
$ws.Range("O5").Value = 0.831988775142814
$ws.Range("P5").Value = 0.831988775142814
$ws.Range("Q5").Value = 401.7780670790266
$ws.Range("R5").Value = 3616.00260371124
$ws.Range("S5").Value = 0.2319092787196376
$ws.Range("T5").Value = 0.2319092787196376

# Row 6
$ws.Range("I6").Value = 0.44716501655323
$ws.Range("J6").Value = 0.4471650165532299
$ws.Range("M6").Value = 0.05237733333333333
$ws.Range("N6").Value = 0.157132
$ws.Range("O6").Value = 0.01249705432598326
$ws.Range("P6").Value = 0.01249705432598326
$ws.Range("Q6").Value = 9.681520677287999
$ws.Range("R6").Value = 87.133686095592
$ws.Range("S6").Value = 0.005588245504544918
$ws.Range("T6").Value = 0.005588245504544917

# Row 7
$ws.Range("I7").Value = 0.44716501655323
$ws.Range("J7").Value = 0.4471650165532299
$ws.Range("O7").Value = 0.1346970137168397
$ws.Range("P7").Value = 0.1346970137168397
$ws.Range("S7").Value = 0.06023179236836129
$ws.Range("T7").Value = 0.06023179236836128

# Row 8
$ws.Range("I8").Value = 0.44716501655323
$ws.Range("J8").Value = 0.4471650165532299
$ws.Range("O8").Value = 0.020817156814363
$ws.Range("P8").Value = 0.020817156814363
$ws.Range("S8").Value = 0.009308704271485817
$ws.Range("T8").Value = 0.009308704271485816

# Row 9
$ws.Range("I9").Value = 0.44716501655323
$ws.Range("J9").Value = 0.4471650165532299
$ws.Range("M9").Value = 3.48701
$ws.Range("N9").Value = 10.46103
$ws.Range("O9").Value = 0.831988775142814
$ws.Range("P9").Value = 0.831988775142814
$ws.Range("Q9").Value = 644.5452119920201
$ws.Range("R9").Value = 5800.90690792818
$ws.Range("S9").Value = 0.372036274408838
$ws.Range("T9").Value = 0.3720362744088379

# Row 10
$ws.Range("G10").Value = 60.55095666666667
$ws.Range("H10").Value = 181.65287
$ws.Range("I10").Value = 0.1464834753134679
$ws.Range("J10").Value = 0.1464834753134678
$ws.Range("M10").Value = 0.05237733333333333
$ws.Range("N10").Value = 0.157132
$ws.Range("O10").Value = 0.01249705432598326
$ws.Range("P10").Value = 0.01249705432598326
$ws.Range("Q10").Value = 3.171497640982222
$ws.Range("R10").Value = 28.54347876884
$ws.Range("S10").Value = 0.001830611948851235
$ws.Range("T10").Value = 0.001830611948851235

# Row 11
$ws.Range("G11").Value = 60.55095666666667
$ws.Range("H11").Value = 181.65287
$ws.Range("I11").Value = 0.1464834753134679
$ws.Range("J11").Value = 0.1464834753134678
$ws.Range("O11").Value = 0.1346970137168397
$ws.Range("P11").Value = 0.1346970137168397
$ws.Range("Q11").Value = 34.18335634199111
$ws.Range("R11").Value = 307.65020707792
$ws.Range("S11").Value = 0.01973088668358854
$ws.Range("T11").Value = 0.01973088668358853

# Row 12
$ws.Range("G12").Value = 60.55095666666667
$ws.Range("H12").Value = 181.65287
$ws.Range("I12").Value = 0.1464834753134679
$ws.Range("J12").Value = 0.1464834753134678
$ws.Range("O12").Value = 0.020817156814363
$ws.Range("P12").Value = 0.020817156814363
$ws.Range("Q12").Value = 5.282970050905556
$ws.Range("R12").Value = 47.54673045815
$ws.Range("S12").Value = 0.003049369476313332
$ws.Range("T12").Value = 0.003049369476313332

# Row 13
$ws.Range("G13").Value = 60.55095666666667
$ws.Range("H13").Value = 181.65287
$ws.Range("I13").Value = 0.1464834753134679
$ws.Range("J13").Value = 0.1464834753134678
$ws.Range("M13").Value = 3.48701
$ws.Range("N13").Value = 10.46103
$ws.Range("O13").Value = 0.831988775142814
$ws.Range("P13").Value = 0.831988775142814
$ws.Range("Q13").Value = 211.1417914062334
$ws.Range("R13").Value = 1900.2761226561
$ws.Range("S13").Value = 0.1218726072047148
$ws.Range("T13").Value = 0.1218726072047147

# Row 14
$ws.Range("G14").Value = 52.74960833333333
$ws.Range("H14").Value = 158.248825
$ws.Range("I14").Value = 0.1276106336788006
$ws.Range("J14").Value = 0.1276106336788006
$ws.Range("M14").Value = 0.05237733333333333
$ws.Range("N14").Value = 0.157132
$ws.Range("O14").Value = 0.01249705432598326
$ws.Range("P14").Value = 0.01249705432598326
$ws.Range("Q14").Value = 2.762883818877778
$ws.Range("R14").Value = 24.8659543699
$ws.Range("S14").Value = 0.00159475702165712
$ws.Range("T14").Value = 0.00159475702165712

# Row 15
$ws.Range("G15").Value = 52.74960833333333
$ws.Range("H15").Value = 158.248825
$ws.Range("I15").Value = 0.1276106336788006
$ws.Range("J15").Value = 0.1276106336788006
$ws.Range("O15").Value = 0.1346970137168397
$ws.Range("P15").Value = 0.1346970137168397
$ws.Range("Q15").Value = 29.77919355568889
$ws.Range("R15").Value = 268.0127420012
$ws.Range("S15").Value = 0.01718877127504802
$ws.Range("T15").Value = 0.01718877127504802

# Row 16
$ws.Range("G16").Value = 52.74960833333333
$ws.Range("H16").Value = 158.248825
$ws.Range("I16").Value = 0.1276106336788006
$ws.Range("J16").Value = 0.1276106336788006
$ws.Range("O16").Value = 0.020817156814363
$ws.Range("P16").Value = 0.020817156814363
$ws.Range("Q16").Value = 4.602315411069444
$ws.Range("R16").Value = 41.420838699625
$ws.Range("S16").Value = 0.002656490572471826
$ws.Range("T16").Value = 0.002656490572471826

# Row 17
$ws.Range("G17").Value = 52.74960833333333
$ws.Range("H17").Value = 158.248825
$ws.Range("I17").Value = 0.1276106336788006
$ws.Range("J17").Value = 0.1276106336788006
$ws.Range("M17").Value = 3.48701
$ws.Range("N17").Value = 10.46103
$ws.Range("O17").Value = 0.831988775142814
$ws.Range("P17").Value = 0.831988775142814
$ws.Range("Q17").Value = 183.9384117544167
$ws.Range("R17").Value = 1655.44570578975
$ws.Range("S17").Value = 0.1061706148096237
$ws.Range("T17").Value = 0.1061706148096237
